# Apply weekly Fruit/Vegetable price-history update for "Femacal de La Calera - Tuna".
# Three new observation rows are inserted at the top of the data block (new rows 61-63),
# pushing the previously existing rows 61-96 down to rows 64-99, and one additional new
# row is appended at the end (row 100, previously-last row 97 data unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61
$ws.Cells.Item(61, 1).Value = 3
$ws.Cells.Item(61, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(61, 3).Value = 'Coquimbo'
$ws.Cells.Item(61, 4).Value = 44606
$ws.Cells.Item(61, 5).Value = 5
$ws.Cells.Item(61, 6).Value = 'Fruta'
$ws.Cells.Item(61, 7).Value = 100107
$ws.Cells.Item(61, 8).Value = 'Otros'
$ws.Cells.Item(61, 9).Value = 100107011
$ws.Cells.Item(61, 10).Value = 'Tuna'
$ws.Cells.Item(61, 11).Value = 'Sin especificar'
$ws.Cells.Item(61, 12).Value = 'Especial'
$ws.Cells.Item(61, 13).Value = 45
$ws.Cells.Item(61, 14).Value = 18000
$ws.Cells.Item(61, 15).Value = 18000
$ws.Cells.Item(61, 16).Value = 18000
$ws.Cells.Item(61, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(61, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(61, 19).Value = 1125
$ws.Cells.Item(61, 20).Value = 16

# Row 62
$ws.Cells.Item(62, 1).Value = 3
$ws.Cells.Item(62, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(62, 3).Value = 'Coquimbo'
$ws.Cells.Item(62, 4).Value = 44606
$ws.Cells.Item(62, 5).Value = 5
$ws.Cells.Item(62, 6).Value = 'Fruta'
$ws.Cells.Item(62, 7).Value = 100107
$ws.Cells.Item(62, 8).Value = 'Otros'
$ws.Cells.Item(62, 9).Value = 100107011
$ws.Cells.Item(62, 10).Value = 'Tuna'
$ws.Cells.Item(62, 11).Value = 'Sin especificar'
$ws.Cells.Item(62, 12).Value = 'Primera'
$ws.Cells.Item(62, 13).Value = 57
$ws.Cells.Item(62, 14).Value = 16000
$ws.Cells.Item(62, 15).Value = 16000
$ws.Cells.Item(62, 16).Value = 16000
$ws.Cells.Item(62, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(62, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(62, 19).Value = 1000
$ws.Cells.Item(62, 20).Value = 16

# Row 63
$ws.Cells.Item(63, 1).Value = 3
$ws.Cells.Item(63, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(63, 3).Value = 'Coquimbo'
$ws.Cells.Item(63, 4).Value = 44606
$ws.Cells.Item(63, 5).Value = 5
$ws.Cells.Item(63, 6).Value = 'Fruta'
$ws.Cells.Item(63, 7).Value = 100107
$ws.Cells.Item(63, 8).Value = 'Otros'
$ws.Cells.Item(63, 9).Value = 100107011
$ws.Cells.Item(63, 10).Value = 'Tuna'
$ws.Cells.Item(63, 11).Value = 'Sin especificar'
$ws.Cells.Item(63, 12).Value = 'Segunda'
$ws.Cells.Item(63, 13).Value = 50
$ws.Cells.Item(63, 14).Value = 14000
$ws.Cells.Item(63, 15).Value = 14000
$ws.Cells.Item(63, 16).Value = 14000
$ws.Cells.Item(63, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(63, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(63, 19).Value = 875
$ws.Cells.Item(63, 20).Value = 16

# Row 64
$ws.Cells.Item(64, 1).Value = 3
$ws.Cells.Item(64, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(64, 3).Value = 'Coquimbo'
$ws.Cells.Item(64, 4).Value = 44238
$ws.Cells.Item(64, 5).Value = 5
$ws.Cells.Item(64, 6).Value = 'Fruta'
$ws.Cells.Item(64, 7).Value = 100107
$ws.Cells.Item(64, 8).Value = 'Otros'
$ws.Cells.Item(64, 9).Value = 100107011
$ws.Cells.Item(64, 10).Value = 'Tuna'
$ws.Cells.Item(64, 11).Value = 'Sin especificar'
$ws.Cells.Item(64, 12).Value = 'Primera'
$ws.Cells.Item(64, 13).Value = 90
$ws.Cells.Item(64, 14).Value = 12000
$ws.Cells.Item(64, 15).Value = 12000
$ws.Cells.Item(64, 16).Value = 12000
$ws.Cells.Item(64, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(64, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(64, 19).Value = 750
$ws.Cells.Item(64, 20).Value = 16

# Row 65
$ws.Cells.Item(65, 1).Value = 3
$ws.Cells.Item(65, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(65, 3).Value = 'Coquimbo'
$ws.Cells.Item(65, 4).Value = 44242
$ws.Cells.Item(65, 5).Value = 5
$ws.Cells.Item(65, 6).Value = 'Fruta'
$ws.Cells.Item(65, 7).Value = 100107
$ws.Cells.Item(65, 8).Value = 'Otros'
$ws.Cells.Item(65, 9).Value = 100107011
$ws.Cells.Item(65, 10).Value = 'Tuna'
$ws.Cells.Item(65, 11).Value = 'Sin especificar'
$ws.Cells.Item(65, 12).Value = 'Primera'
$ws.Cells.Item(65, 13).Value = 68
$ws.Cells.Item(65, 14).Value = 12000
$ws.Cells.Item(65, 15).Value = 12000
$ws.Cells.Item(65, 16).Value = 12000
$ws.Cells.Item(65, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(65, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(65, 19).Value = 750
$ws.Cells.Item(65, 20).Value = 16

# Row 66
$ws.Cells.Item(66, 1).Value = 3
$ws.Cells.Item(66, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(66, 3).Value = 'Coquimbo'
$ws.Cells.Item(66, 4).Value = 44249
$ws.Cells.Item(66, 5).Value = 5
$ws.Cells.Item(66, 6).Value = 'Fruta'
$ws.Cells.Item(66, 7).Value = 100107
$ws.Cells.Item(66, 8).Value = 'Otros'
$ws.Cells.Item(66, 9).Value = 100107011
$ws.Cells.Item(66, 10).Value = 'Tuna'
$ws.Cells.Item(66, 11).Value = 'Sin especificar'
$ws.Cells.Item(66, 12).Value = 'Primera'
$ws.Cells.Item(66, 13).Value = 80
$ws.Cells.Item(66, 14).Value = 12000
$ws.Cells.Item(66, 15).Value = 12000
$ws.Cells.Item(66, 16).Value = 12000
$ws.Cells.Item(66, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(66, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(66, 19).Value = 750
$ws.Cells.Item(66, 20).Value = 16

# Row 67
$ws.Cells.Item(67, 1).Value = 3
$ws.Cells.Item(67, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(67, 3).Value = 'Coquimbo'
$ws.Cells.Item(67, 4).Value = 44252
$ws.Cells.Item(67, 5).Value = 5
$ws.Cells.Item(67, 6).Value = 'Fruta'
$ws.Cells.Item(67, 7).Value = 100107
$ws.Cells.Item(67, 8).Value = 'Otros'
$ws.Cells.Item(67, 9).Value = 100107011
$ws.Cells.Item(67, 10).Value = 'Tuna'
$ws.Cells.Item(67, 11).Value = 'Sin especificar'
$ws.Cells.Item(67, 12).Value = 'Primera'
$ws.Cells.Item(67, 13).Value = 87
$ws.Cells.Item(67, 14).Value = 15000
$ws.Cells.Item(67, 15).Value = 15000
$ws.Cells.Item(67, 16).Value = 15000
$ws.Cells.Item(67, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(67, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(67, 19).Value = 938
$ws.Cells.Item(67, 20).Value = 16

# Row 68
$ws.Cells.Item(68, 1).Value = 3
$ws.Cells.Item(68, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(68, 3).Value = 'Coquimbo'
$ws.Cells.Item(68, 4).Value = 44271
$ws.Cells.Item(68, 5).Value = 5
$ws.Cells.Item(68, 6).Value = 'Fruta'
$ws.Cells.Item(68, 7).Value = 100107
$ws.Cells.Item(68, 8).Value = 'Otros'
$ws.Cells.Item(68, 9).Value = 100107011
$ws.Cells.Item(68, 10).Value = 'Tuna'
$ws.Cells.Item(68, 11).Value = 'Sin especificar'
$ws.Cells.Item(68, 12).Value = 'Especial'
$ws.Cells.Item(68, 13).Value = 50
$ws.Cells.Item(68, 14).Value = 12000
$ws.Cells.Item(68, 15).Value = 12000
$ws.Cells.Item(68, 16).Value = 12000
$ws.Cells.Item(68, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(68, 18).Value = 'Cabildo'
$ws.Cells.Item(68, 19).Value = 750
$ws.Cells.Item(68, 20).Value = 16

# Row 69
$ws.Cells.Item(69, 1).Value = 3
$ws.Cells.Item(69, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(69, 3).Value = 'Coquimbo'
$ws.Cells.Item(69, 4).Value = 44271
$ws.Cells.Item(69, 5).Value = 5
$ws.Cells.Item(69, 6).Value = 'Fruta'
$ws.Cells.Item(69, 7).Value = 100107
$ws.Cells.Item(69, 8).Value = 'Otros'
$ws.Cells.Item(69, 9).Value = 100107011
$ws.Cells.Item(69, 10).Value = 'Tuna'
$ws.Cells.Item(69, 11).Value = 'Sin especificar'
$ws.Cells.Item(69, 12).Value = 'Primera'
$ws.Cells.Item(69, 13).Value = 54
$ws.Cells.Item(69, 14).Value = 10000
$ws.Cells.Item(69, 15).Value = 10000
$ws.Cells.Item(69, 16).Value = 10000
$ws.Cells.Item(69, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(69, 18).Value = 'Cabildo'
$ws.Cells.Item(69, 19).Value = 625
$ws.Cells.Item(69, 20).Value = 16

# Row 70
$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(70, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(70, 3).Value = 'Coquimbo'
$ws.Cells.Item(70, 4).Value = 44271
$ws.Cells.Item(70, 5).Value = 5
$ws.Cells.Item(70, 6).Value = 'Fruta'
$ws.Cells.Item(70, 7).Value = 100107
$ws.Cells.Item(70, 8).Value = 'Otros'
$ws.Cells.Item(70, 9).Value = 100107011
$ws.Cells.Item(70, 10).Value = 'Tuna'
$ws.Cells.Item(70, 11).Value = 'Sin especificar'
$ws.Cells.Item(70, 12).Value = 'Segunda'
$ws.Cells.Item(70, 13).Value = 48
$ws.Cells.Item(70, 14).Value = 8000
$ws.Cells.Item(70, 15).Value = 8000
$ws.Cells.Item(70, 16).Value = 8000
$ws.Cells.Item(70, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(70, 18).Value = 'Cabildo'
$ws.Cells.Item(70, 19).Value = 500
$ws.Cells.Item(70, 20).Value = 16

# Row 71
$ws.Cells.Item(71, 1).Value = 3
$ws.Cells.Item(71, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(71, 3).Value = 'Coquimbo'
$ws.Cells.Item(71, 4).Value = 44315
$ws.Cells.Item(71, 5).Value = 5
$ws.Cells.Item(71, 6).Value = 'Fruta'
$ws.Cells.Item(71, 7).Value = 100107
$ws.Cells.Item(71, 8).Value = 'Otros'
$ws.Cells.Item(71, 9).Value = 100107011
$ws.Cells.Item(71, 10).Value = 'Tuna'
$ws.Cells.Item(71, 11).Value = 'Sin especificar'
$ws.Cells.Item(71, 12).Value = 'Especial'
$ws.Cells.Item(71, 13).Value = 70
$ws.Cells.Item(71, 14).Value = 18000
$ws.Cells.Item(71, 15).Value = 18000
$ws.Cells.Item(71, 16).Value = 18000
$ws.Cells.Item(71, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(71, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(71, 19).Value = 1125
$ws.Cells.Item(71, 20).Value = 16

# Row 72
$ws.Cells.Item(72, 1).Value = 3
$ws.Cells.Item(72, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(72, 3).Value = 'Coquimbo'
$ws.Cells.Item(72, 4).Value = 44315
$ws.Cells.Item(72, 5).Value = 5
$ws.Cells.Item(72, 6).Value = 'Fruta'
$ws.Cells.Item(72, 7).Value = 100107
$ws.Cells.Item(72, 8).Value = 'Otros'
$ws.Cells.Item(72, 9).Value = 100107011
$ws.Cells.Item(72, 10).Value = 'Tuna'
$ws.Cells.Item(72, 11).Value = 'Sin especificar'
$ws.Cells.Item(72, 12).Value = 'Primera'
$ws.Cells.Item(72, 13).Value = 75
$ws.Cells.Item(72, 14).Value = 15000
$ws.Cells.Item(72, 15).Value = 15000
$ws.Cells.Item(72, 16).Value = 15000
$ws.Cells.Item(72, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(72, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(72, 19).Value = 938
$ws.Cells.Item(72, 20).Value = 16

# Row 73
$ws.Cells.Item(73, 1).Value = 3
$ws.Cells.Item(73, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(73, 3).Value = 'Coquimbo'
$ws.Cells.Item(73, 4).Value = 44315
$ws.Cells.Item(73, 5).Value = 5
$ws.Cells.Item(73, 6).Value = 'Fruta'
$ws.Cells.Item(73, 7).Value = 100107
$ws.Cells.Item(73, 8).Value = 'Otros'
$ws.Cells.Item(73, 9).Value = 100107011
$ws.Cells.Item(73, 10).Value = 'Tuna'
$ws.Cells.Item(73, 11).Value = 'Sin especificar'
$ws.Cells.Item(73, 12).Value = 'Segunda'
$ws.Cells.Item(73, 13).Value = 68
$ws.Cells.Item(73, 14).Value = 12000
$ws.Cells.Item(73, 15).Value = 12000
$ws.Cells.Item(73, 16).Value = 12000
$ws.Cells.Item(73, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(73, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(73, 19).Value = 750
$ws.Cells.Item(73, 20).Value = 16

# Row 74
$ws.Cells.Item(74, 1).Value = 3
$ws.Cells.Item(74, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(74, 3).Value = 'Coquimbo'
$ws.Cells.Item(74, 4).Value = 44280
$ws.Cells.Item(74, 5).Value = 5
$ws.Cells.Item(74, 6).Value = 'Fruta'
$ws.Cells.Item(74, 7).Value = 100107
$ws.Cells.Item(74, 8).Value = 'Otros'
$ws.Cells.Item(74, 9).Value = 100107011
$ws.Cells.Item(74, 10).Value = 'Tuna'
$ws.Cells.Item(74, 11).Value = 'Sin especificar'
$ws.Cells.Item(74, 12).Value = 'Especial'
$ws.Cells.Item(74, 13).Value = 65
$ws.Cells.Item(74, 14).Value = 12000
$ws.Cells.Item(74, 15).Value = 12000
$ws.Cells.Item(74, 16).Value = 12000
$ws.Cells.Item(74, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(74, 18).Value = 'Cabildo'
$ws.Cells.Item(74, 19).Value = 750
$ws.Cells.Item(74, 20).Value = 16

# Row 75
$ws.Cells.Item(75, 1).Value = 3
$ws.Cells.Item(75, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(75, 3).Value = 'Coquimbo'
$ws.Cells.Item(75, 4).Value = 44280
$ws.Cells.Item(75, 5).Value = 5
$ws.Cells.Item(75, 6).Value = 'Fruta'
$ws.Cells.Item(75, 7).Value = 100107
$ws.Cells.Item(75, 8).Value = 'Otros'
$ws.Cells.Item(75, 9).Value = 100107011
$ws.Cells.Item(75, 10).Value = 'Tuna'
$ws.Cells.Item(75, 11).Value = 'Sin especificar'
$ws.Cells.Item(75, 12).Value = 'Primera'
$ws.Cells.Item(75, 13).Value = 70
$ws.Cells.Item(75, 14).Value = 10000
$ws.Cells.Item(75, 15).Value = 10000
$ws.Cells.Item(75, 16).Value = 10000
$ws.Cells.Item(75, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(75, 18).Value = 'Cabildo'
$ws.Cells.Item(75, 19).Value = 625
$ws.Cells.Item(75, 20).Value = 16

# Row 76
$ws.Cells.Item(76, 1).Value = 3
$ws.Cells.Item(76, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(76, 3).Value = 'Coquimbo'
$ws.Cells.Item(76, 4).Value = 44280
$ws.Cells.Item(76, 5).Value = 5
$ws.Cells.Item(76, 6).Value = 'Fruta'
$ws.Cells.Item(76, 7).Value = 100107
$ws.Cells.Item(76, 8).Value = 'Otros'
$ws.Cells.Item(76, 9).Value = 100107011
$ws.Cells.Item(76, 10).Value = 'Tuna'
$ws.Cells.Item(76, 11).Value = 'Sin especificar'
$ws.Cells.Item(76, 12).Value = 'Segunda'
$ws.Cells.Item(76, 13).Value = 68
$ws.Cells.Item(76, 14).Value = 8000
$ws.Cells.Item(76, 15).Value = 8000
$ws.Cells.Item(76, 16).Value = 8000
$ws.Cells.Item(76, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(76, 18).Value = 'Cabildo'
$ws.Cells.Item(76, 19).Value = 500
$ws.Cells.Item(76, 20).Value = 16

# Row 77
$ws.Cells.Item(77, 1).Value = 3
$ws.Cells.Item(77, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(77, 3).Value = 'Coquimbo'
$ws.Cells.Item(77, 4).Value = 44581
$ws.Cells.Item(77, 5).Value = 5
$ws.Cells.Item(77, 6).Value = 'Fruta'
$ws.Cells.Item(77, 7).Value = 100107
$ws.Cells.Item(77, 8).Value = 'Otros'
$ws.Cells.Item(77, 9).Value = 100107011
$ws.Cells.Item(77, 10).Value = 'Tuna'
$ws.Cells.Item(77, 11).Value = 'Sin especificar'
$ws.Cells.Item(77, 12).Value = 'Primera'
$ws.Cells.Item(77, 13).Value = 50
$ws.Cells.Item(77, 14).Value = 20000
$ws.Cells.Item(77, 15).Value = 20000
$ws.Cells.Item(77, 16).Value = 20000
$ws.Cells.Item(77, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(77, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(77, 19).Value = 1000
$ws.Cells.Item(77, 20).Value = 20

# Row 78
$ws.Cells.Item(78, 1).Value = 3
$ws.Cells.Item(78, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(78, 3).Value = 'Coquimbo'
$ws.Cells.Item(78, 4).Value = 44294
$ws.Cells.Item(78, 5).Value = 5
$ws.Cells.Item(78, 6).Value = 'Fruta'
$ws.Cells.Item(78, 7).Value = 100107
$ws.Cells.Item(78, 8).Value = 'Otros'
$ws.Cells.Item(78, 9).Value = 100107011
$ws.Cells.Item(78, 10).Value = 'Tuna'
$ws.Cells.Item(78, 11).Value = 'Sin especificar'
$ws.Cells.Item(78, 12).Value = 'Especial'
$ws.Cells.Item(78, 13).Value = 65
$ws.Cells.Item(78, 14).Value = 18000
$ws.Cells.Item(78, 15).Value = 18000
$ws.Cells.Item(78, 16).Value = 18000
$ws.Cells.Item(78, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(78, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(78, 19).Value = 1125
$ws.Cells.Item(78, 20).Value = 16

# Row 79
$ws.Cells.Item(79, 1).Value = 3
$ws.Cells.Item(79, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(79, 3).Value = 'Coquimbo'
$ws.Cells.Item(79, 4).Value = 44294
$ws.Cells.Item(79, 5).Value = 5
$ws.Cells.Item(79, 6).Value = 'Fruta'
$ws.Cells.Item(79, 7).Value = 100107
$ws.Cells.Item(79, 8).Value = 'Otros'
$ws.Cells.Item(79, 9).Value = 100107011
$ws.Cells.Item(79, 10).Value = 'Tuna'
$ws.Cells.Item(79, 11).Value = 'Sin especificar'
$ws.Cells.Item(79, 12).Value = 'Primera'
$ws.Cells.Item(79, 13).Value = 67
$ws.Cells.Item(79, 14).Value = 15000
$ws.Cells.Item(79, 15).Value = 15000
$ws.Cells.Item(79, 16).Value = 15000
$ws.Cells.Item(79, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(79, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(79, 19).Value = 938
$ws.Cells.Item(79, 20).Value = 16

# Row 80
$ws.Cells.Item(80, 1).Value = 3
$ws.Cells.Item(80, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(80, 3).Value = 'Coquimbo'
$ws.Cells.Item(80, 4).Value = 44294
$ws.Cells.Item(80, 5).Value = 5
$ws.Cells.Item(80, 6).Value = 'Fruta'
$ws.Cells.Item(80, 7).Value = 100107
$ws.Cells.Item(80, 8).Value = 'Otros'
$ws.Cells.Item(80, 9).Value = 100107011
$ws.Cells.Item(80, 10).Value = 'Tuna'
$ws.Cells.Item(80, 11).Value = 'Sin especificar'
$ws.Cells.Item(80, 12).Value = 'Segunda'
$ws.Cells.Item(80, 13).Value = 60
$ws.Cells.Item(80, 14).Value = 12000
$ws.Cells.Item(80, 15).Value = 12000
$ws.Cells.Item(80, 16).Value = 12000
$ws.Cells.Item(80, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(80, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(80, 19).Value = 750
$ws.Cells.Item(80, 20).Value = 16

# Row 81
$ws.Cells.Item(81, 1).Value = 3
$ws.Cells.Item(81, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(81, 3).Value = 'Coquimbo'
$ws.Cells.Item(81, 4).Value = 44232
$ws.Cells.Item(81, 5).Value = 5
$ws.Cells.Item(81, 6).Value = 'Fruta'
$ws.Cells.Item(81, 7).Value = 100107
$ws.Cells.Item(81, 8).Value = 'Otros'
$ws.Cells.Item(81, 9).Value = 100107011
$ws.Cells.Item(81, 10).Value = 'Tuna'
$ws.Cells.Item(81, 11).Value = 'Sin especificar'
$ws.Cells.Item(81, 12).Value = 'Primera'
$ws.Cells.Item(81, 13).Value = 70
$ws.Cells.Item(81, 14).Value = 14000
$ws.Cells.Item(81, 15).Value = 14000
$ws.Cells.Item(81, 16).Value = 14000
$ws.Cells.Item(81, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(81, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(81, 19).Value = 875
$ws.Cells.Item(81, 20).Value = 16

# Row 82
$ws.Cells.Item(82, 1).Value = 3
$ws.Cells.Item(82, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(82, 3).Value = 'Coquimbo'
$ws.Cells.Item(82, 4).Value = 44232
$ws.Cells.Item(82, 5).Value = 5
$ws.Cells.Item(82, 6).Value = 'Fruta'
$ws.Cells.Item(82, 7).Value = 100107
$ws.Cells.Item(82, 8).Value = 'Otros'
$ws.Cells.Item(82, 9).Value = 100107011
$ws.Cells.Item(82, 10).Value = 'Tuna'
$ws.Cells.Item(82, 11).Value = 'Sin especificar'
$ws.Cells.Item(82, 12).Value = 'Segunda'
$ws.Cells.Item(82, 13).Value = 120
$ws.Cells.Item(82, 14).Value = 12000
$ws.Cells.Item(82, 15).Value = 12000
$ws.Cells.Item(82, 16).Value = 12000
$ws.Cells.Item(82, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(82, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(82, 19).Value = 750
$ws.Cells.Item(82, 20).Value = 16

# Row 83
$ws.Cells.Item(83, 1).Value = 3
$ws.Cells.Item(83, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(83, 3).Value = 'Coquimbo'
$ws.Cells.Item(83, 4).Value = 44279
$ws.Cells.Item(83, 5).Value = 5
$ws.Cells.Item(83, 6).Value = 'Fruta'
$ws.Cells.Item(83, 7).Value = 100107
$ws.Cells.Item(83, 8).Value = 'Otros'
$ws.Cells.Item(83, 9).Value = 100107011
$ws.Cells.Item(83, 10).Value = 'Tuna'
$ws.Cells.Item(83, 11).Value = 'Sin especificar'
$ws.Cells.Item(83, 12).Value = 'Especial'
$ws.Cells.Item(83, 13).Value = 65
$ws.Cells.Item(83, 14).Value = 12000
$ws.Cells.Item(83, 15).Value = 12000
$ws.Cells.Item(83, 16).Value = 12000
$ws.Cells.Item(83, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(83, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(83, 19).Value = 750
$ws.Cells.Item(83, 20).Value = 16

# Row 84
$ws.Cells.Item(84, 1).Value = 3
$ws.Cells.Item(84, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(84, 3).Value = 'Coquimbo'
$ws.Cells.Item(84, 4).Value = 44279
$ws.Cells.Item(84, 5).Value = 5
$ws.Cells.Item(84, 6).Value = 'Fruta'
$ws.Cells.Item(84, 7).Value = 100107
$ws.Cells.Item(84, 8).Value = 'Otros'
$ws.Cells.Item(84, 9).Value = 100107011
$ws.Cells.Item(84, 10).Value = 'Tuna'
$ws.Cells.Item(84, 11).Value = 'Sin especificar'
$ws.Cells.Item(84, 12).Value = 'Primera'
$ws.Cells.Item(84, 13).Value = 78
$ws.Cells.Item(84, 14).Value = 10000
$ws.Cells.Item(84, 15).Value = 10000
$ws.Cells.Item(84, 16).Value = 10000
$ws.Cells.Item(84, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(84, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(84, 19).Value = 625
$ws.Cells.Item(84, 20).Value = 16

# Row 85
$ws.Cells.Item(85, 1).Value = 3
$ws.Cells.Item(85, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(85, 3).Value = 'Coquimbo'
$ws.Cells.Item(85, 4).Value = 44279
$ws.Cells.Item(85, 5).Value = 5
$ws.Cells.Item(85, 6).Value = 'Fruta'
$ws.Cells.Item(85, 7).Value = 100107
$ws.Cells.Item(85, 8).Value = 'Otros'
$ws.Cells.Item(85, 9).Value = 100107011
$ws.Cells.Item(85, 10).Value = 'Tuna'
$ws.Cells.Item(85, 11).Value = 'Sin especificar'
$ws.Cells.Item(85, 12).Value = 'Segunda'
$ws.Cells.Item(85, 13).Value = 70
$ws.Cells.Item(85, 14).Value = 8000
$ws.Cells.Item(85, 15).Value = 8000
$ws.Cells.Item(85, 16).Value = 8000
$ws.Cells.Item(85, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(85, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(85, 19).Value = 500
$ws.Cells.Item(85, 20).Value = 16

# Row 86
$ws.Cells.Item(86, 1).Value = 3
$ws.Cells.Item(86, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(86, 3).Value = 'Coquimbo'
$ws.Cells.Item(86, 4).Value = 44301
$ws.Cells.Item(86, 5).Value = 5
$ws.Cells.Item(86, 6).Value = 'Fruta'
$ws.Cells.Item(86, 7).Value = 100107
$ws.Cells.Item(86, 8).Value = 'Otros'
$ws.Cells.Item(86, 9).Value = 100107011
$ws.Cells.Item(86, 10).Value = 'Tuna'
$ws.Cells.Item(86, 11).Value = 'Sin especificar'
$ws.Cells.Item(86, 12).Value = 'Especial'
$ws.Cells.Item(86, 13).Value = 56
$ws.Cells.Item(86, 14).Value = 18000
$ws.Cells.Item(86, 15).Value = 18000
$ws.Cells.Item(86, 16).Value = 18000
$ws.Cells.Item(86, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(86, 18).Value = 'Cabildo'
$ws.Cells.Item(86, 19).Value = 1125
$ws.Cells.Item(86, 20).Value = 16

# Row 87
$ws.Cells.Item(87, 1).Value = 3
$ws.Cells.Item(87, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(87, 3).Value = 'Coquimbo'
$ws.Cells.Item(87, 4).Value = 44301
$ws.Cells.Item(87, 5).Value = 5
$ws.Cells.Item(87, 6).Value = 'Fruta'
$ws.Cells.Item(87, 7).Value = 100107
$ws.Cells.Item(87, 8).Value = 'Otros'
$ws.Cells.Item(87, 9).Value = 100107011
$ws.Cells.Item(87, 10).Value = 'Tuna'
$ws.Cells.Item(87, 11).Value = 'Sin especificar'
$ws.Cells.Item(87, 12).Value = 'Primera'
$ws.Cells.Item(87, 13).Value = 68
$ws.Cells.Item(87, 14).Value = 15000
$ws.Cells.Item(87, 15).Value = 15000
$ws.Cells.Item(87, 16).Value = 15000
$ws.Cells.Item(87, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(87, 18).Value = 'Cabildo'
$ws.Cells.Item(87, 19).Value = 938
$ws.Cells.Item(87, 20).Value = 16

# Row 88
$ws.Cells.Item(88, 1).Value = 3
$ws.Cells.Item(88, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(88, 3).Value = 'Coquimbo'
$ws.Cells.Item(88, 4).Value = 44301
$ws.Cells.Item(88, 5).Value = 5
$ws.Cells.Item(88, 6).Value = 'Fruta'
$ws.Cells.Item(88, 7).Value = 100107
$ws.Cells.Item(88, 8).Value = 'Otros'
$ws.Cells.Item(88, 9).Value = 100107011
$ws.Cells.Item(88, 10).Value = 'Tuna'
$ws.Cells.Item(88, 11).Value = 'Sin especificar'
$ws.Cells.Item(88, 12).Value = 'Segunda'
$ws.Cells.Item(88, 13).Value = 60
$ws.Cells.Item(88, 14).Value = 12000
$ws.Cells.Item(88, 15).Value = 12000
$ws.Cells.Item(88, 16).Value = 12000
$ws.Cells.Item(88, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(88, 18).Value = 'Cabildo'
$ws.Cells.Item(88, 19).Value = 750
$ws.Cells.Item(88, 20).Value = 16

# Row 89
$ws.Cells.Item(89, 1).Value = 3
$ws.Cells.Item(89, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(89, 3).Value = 'Coquimbo'
$ws.Cells.Item(89, 4).Value = 44221
$ws.Cells.Item(89, 5).Value = 5
$ws.Cells.Item(89, 6).Value = 'Fruta'
$ws.Cells.Item(89, 7).Value = 100107
$ws.Cells.Item(89, 8).Value = 'Otros'
$ws.Cells.Item(89, 9).Value = 100107011
$ws.Cells.Item(89, 10).Value = 'Tuna'
$ws.Cells.Item(89, 11).Value = 'Sin especificar'
$ws.Cells.Item(89, 12).Value = 'Primera'
$ws.Cells.Item(89, 13).Value = 68
$ws.Cells.Item(89, 14).Value = 15000
$ws.Cells.Item(89, 15).Value = 15000
$ws.Cells.Item(89, 16).Value = 15000
$ws.Cells.Item(89, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(89, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(89, 19).Value = 938
$ws.Cells.Item(89, 20).Value = 16

# Row 90
$ws.Cells.Item(90, 1).Value = 3
$ws.Cells.Item(90, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(90, 3).Value = 'Coquimbo'
$ws.Cells.Item(90, 4).Value = 44221
$ws.Cells.Item(90, 5).Value = 5
$ws.Cells.Item(90, 6).Value = 'Fruta'
$ws.Cells.Item(90, 7).Value = 100107
$ws.Cells.Item(90, 8).Value = 'Otros'
$ws.Cells.Item(90, 9).Value = 100107011
$ws.Cells.Item(90, 10).Value = 'Tuna'
$ws.Cells.Item(90, 11).Value = 'Sin especificar'
$ws.Cells.Item(90, 12).Value = 'Primera'
$ws.Cells.Item(90, 13).Value = 68
$ws.Cells.Item(90, 14).Value = 18000
$ws.Cells.Item(90, 15).Value = 18000
$ws.Cells.Item(90, 16).Value = 18000
$ws.Cells.Item(90, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(90, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(90, 19).Value = 1000
$ws.Cells.Item(90, 20).Value = 18

# Row 91
$ws.Cells.Item(91, 1).Value = 3
$ws.Cells.Item(91, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(91, 3).Value = 'Coquimbo'
$ws.Cells.Item(91, 4).Value = 44272
$ws.Cells.Item(91, 5).Value = 5
$ws.Cells.Item(91, 6).Value = 'Fruta'
$ws.Cells.Item(91, 7).Value = 100107
$ws.Cells.Item(91, 8).Value = 'Otros'
$ws.Cells.Item(91, 9).Value = 100107011
$ws.Cells.Item(91, 10).Value = 'Tuna'
$ws.Cells.Item(91, 11).Value = 'Sin especificar'
$ws.Cells.Item(91, 12).Value = 'Especial'
$ws.Cells.Item(91, 13).Value = 60
$ws.Cells.Item(91, 14).Value = 12000
$ws.Cells.Item(91, 15).Value = 12000
$ws.Cells.Item(91, 16).Value = 12000
$ws.Cells.Item(91, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(91, 18).Value = 'Cabildo'
$ws.Cells.Item(91, 19).Value = 750
$ws.Cells.Item(91, 20).Value = 16

# Row 92
$ws.Cells.Item(92, 1).Value = 3
$ws.Cells.Item(92, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(92, 3).Value = 'Coquimbo'
$ws.Cells.Item(92, 4).Value = 44272
$ws.Cells.Item(92, 5).Value = 5
$ws.Cells.Item(92, 6).Value = 'Fruta'
$ws.Cells.Item(92, 7).Value = 100107
$ws.Cells.Item(92, 8).Value = 'Otros'
$ws.Cells.Item(92, 9).Value = 100107011
$ws.Cells.Item(92, 10).Value = 'Tuna'
$ws.Cells.Item(92, 11).Value = 'Sin especificar'
$ws.Cells.Item(92, 12).Value = 'Primera'
$ws.Cells.Item(92, 13).Value = 145
$ws.Cells.Item(92, 14).Value = 8000
$ws.Cells.Item(92, 15).Value = 10000
$ws.Cells.Item(92, 16).Value = 9034
$ws.Cells.Item(92, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(92, 18).Value = 'Cabildo'
$ws.Cells.Item(92, 19).Value = 565
$ws.Cells.Item(92, 20).Value = 16

# Row 93
$ws.Cells.Item(93, 1).Value = 3
$ws.Cells.Item(93, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(93, 3).Value = 'Coquimbo'
$ws.Cells.Item(93, 4).Value = 44277
$ws.Cells.Item(93, 5).Value = 5
$ws.Cells.Item(93, 6).Value = 'Fruta'
$ws.Cells.Item(93, 7).Value = 100107
$ws.Cells.Item(93, 8).Value = 'Otros'
$ws.Cells.Item(93, 9).Value = 100107011
$ws.Cells.Item(93, 10).Value = 'Tuna'
$ws.Cells.Item(93, 11).Value = 'Sin especificar'
$ws.Cells.Item(93, 12).Value = 'Especial'
$ws.Cells.Item(93, 13).Value = 65
$ws.Cells.Item(93, 14).Value = 13000
$ws.Cells.Item(93, 15).Value = 13000
$ws.Cells.Item(93, 16).Value = 13000
$ws.Cells.Item(93, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(93, 18).Value = 'Cabildo'
$ws.Cells.Item(93, 19).Value = 812
$ws.Cells.Item(93, 20).Value = 16

# Row 94
$ws.Cells.Item(94, 1).Value = 3
$ws.Cells.Item(94, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(94, 3).Value = 'Coquimbo'
$ws.Cells.Item(94, 4).Value = 44277
$ws.Cells.Item(94, 5).Value = 5
$ws.Cells.Item(94, 6).Value = 'Fruta'
$ws.Cells.Item(94, 7).Value = 100107
$ws.Cells.Item(94, 8).Value = 'Otros'
$ws.Cells.Item(94, 9).Value = 100107011
$ws.Cells.Item(94, 10).Value = 'Tuna'
$ws.Cells.Item(94, 11).Value = 'Sin especificar'
$ws.Cells.Item(94, 12).Value = 'Primera'
$ws.Cells.Item(94, 13).Value = 75
$ws.Cells.Item(94, 14).Value = 11000
$ws.Cells.Item(94, 15).Value = 11000
$ws.Cells.Item(94, 16).Value = 11000
$ws.Cells.Item(94, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(94, 18).Value = 'Cabildo'
$ws.Cells.Item(94, 19).Value = 688
$ws.Cells.Item(94, 20).Value = 16

# Row 95
$ws.Cells.Item(95, 1).Value = 3
$ws.Cells.Item(95, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(95, 3).Value = 'Coquimbo'
$ws.Cells.Item(95, 4).Value = 44277
$ws.Cells.Item(95, 5).Value = 5
$ws.Cells.Item(95, 6).Value = 'Fruta'
$ws.Cells.Item(95, 7).Value = 100107
$ws.Cells.Item(95, 8).Value = 'Otros'
$ws.Cells.Item(95, 9).Value = 100107011
$ws.Cells.Item(95, 10).Value = 'Tuna'
$ws.Cells.Item(95, 11).Value = 'Sin especificar'
$ws.Cells.Item(95, 12).Value = 'Segunda'
$ws.Cells.Item(95, 13).Value = 70
$ws.Cells.Item(95, 14).Value = 9000
$ws.Cells.Item(95, 15).Value = 9000
$ws.Cells.Item(95, 16).Value = 9000
$ws.Cells.Item(95, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(95, 18).Value = 'Cabildo'
$ws.Cells.Item(95, 19).Value = 562
$ws.Cells.Item(95, 20).Value = 16

# Row 96
$ws.Cells.Item(96, 1).Value = 3
$ws.Cells.Item(96, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(96, 3).Value = 'Coquimbo'
$ws.Cells.Item(96, 4).Value = 44258
$ws.Cells.Item(96, 5).Value = 5
$ws.Cells.Item(96, 6).Value = 'Fruta'
$ws.Cells.Item(96, 7).Value = 100107
$ws.Cells.Item(96, 8).Value = 'Otros'
$ws.Cells.Item(96, 9).Value = 100107011
$ws.Cells.Item(96, 10).Value = 'Tuna'
$ws.Cells.Item(96, 11).Value = 'Sin especificar'
$ws.Cells.Item(96, 12).Value = 'Primera'
$ws.Cells.Item(96, 13).Value = 60
$ws.Cells.Item(96, 14).Value = 12000
$ws.Cells.Item(96, 15).Value = 12000
$ws.Cells.Item(96, 16).Value = 12000
$ws.Cells.Item(96, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(96, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(96, 19).Value = 750
$ws.Cells.Item(96, 20).Value = 16

# Row 97
$ws.Cells.Item(97, 1).Value = 3
$ws.Cells.Item(97, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(97, 3).Value = 'Coquimbo'
$ws.Cells.Item(97, 4).Value = 44285
$ws.Cells.Item(97, 5).Value = 5
$ws.Cells.Item(97, 6).Value = 'Fruta'
$ws.Cells.Item(97, 7).Value = 100107
$ws.Cells.Item(97, 8).Value = 'Otros'
$ws.Cells.Item(97, 9).Value = 100107011
$ws.Cells.Item(97, 10).Value = 'Tuna'
$ws.Cells.Item(97, 11).Value = 'Sin especificar'
$ws.Cells.Item(97, 12).Value = 'Especial'
$ws.Cells.Item(97, 13).Value = 58
$ws.Cells.Item(97, 14).Value = 12000
$ws.Cells.Item(97, 15).Value = 12000
$ws.Cells.Item(97, 16).Value = 12000
$ws.Cells.Item(97, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(97, 18).Value = 'Cabildo'
$ws.Cells.Item(97, 19).Value = 750
$ws.Cells.Item(97, 20).Value = 16

# Row 98
$ws.Cells.Item(98, 1).Value = 3
$ws.Cells.Item(98, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(98, 3).Value = 'Coquimbo'
$ws.Cells.Item(98, 4).Value = 44285
$ws.Cells.Item(98, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(98, 5).Value = 5
$ws.Cells.Item(98, 6).Value = 'Fruta'
$ws.Cells.Item(98, 7).Value = 100107
$ws.Cells.Item(98, 8).Value = 'Otros'
$ws.Cells.Item(98, 9).Value = 100107011
$ws.Cells.Item(98, 10).Value = 'Tuna'
$ws.Cells.Item(98, 11).Value = 'Sin especificar'
$ws.Cells.Item(98, 12).Value = 'Primera'
$ws.Cells.Item(98, 13).Value = 60
$ws.Cells.Item(98, 14).Value = 10000
$ws.Cells.Item(98, 15).Value = 10000
$ws.Cells.Item(98, 16).Value = 10000
$ws.Cells.Item(98, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(98, 18).Value = 'Cabildo'
$ws.Cells.Item(98, 19).Value = 625
$ws.Cells.Item(98, 20).Value = 16

# Row 99
$ws.Cells.Item(99, 1).Value = 3
$ws.Cells.Item(99, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(99, 3).Value = 'Coquimbo'
$ws.Cells.Item(99, 4).Value = 44285
$ws.Cells.Item(99, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(99, 5).Value = 5
$ws.Cells.Item(99, 6).Value = 'Fruta'
$ws.Cells.Item(99, 7).Value = 100107
$ws.Cells.Item(99, 8).Value = 'Otros'
$ws.Cells.Item(99, 9).Value = 100107011
$ws.Cells.Item(99, 10).Value = 'Tuna'
$ws.Cells.Item(99, 11).Value = 'Sin especificar'
$ws.Cells.Item(99, 12).Value = 'Segunda'
$ws.Cells.Item(99, 13).Value = 65
$ws.Cells.Item(99, 14).Value = 8000
$ws.Cells.Item(99, 15).Value = 8000
$ws.Cells.Item(99, 16).Value = 8000
$ws.Cells.Item(99, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(99, 18).Value = 'Cabildo'
$ws.Cells.Item(99, 19).Value = 500
$ws.Cells.Item(99, 20).Value = 16

# Row 100
$ws.Cells.Item(100, 1).Value = 3
$ws.Cells.Item(100, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(100, 3).Value = 'Coquimbo'
$ws.Cells.Item(100, 4).Value = 44595
$ws.Cells.Item(100, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(100, 5).Value = 5
$ws.Cells.Item(100, 6).Value = 'Fruta'
$ws.Cells.Item(100, 7).Value = 100107
$ws.Cells.Item(100, 8).Value = 'Otros'
$ws.Cells.Item(100, 9).Value = 100107011
$ws.Cells.Item(100, 10).Value = 'Tuna'
$ws.Cells.Item(100, 11).Value = 'Sin especificar'
$ws.Cells.Item(100, 12).Value = 'Primera'
$ws.Cells.Item(100, 13).Value = 50
$ws.Cells.Item(100, 14).Value = 20000
$ws.Cells.Item(100, 15).Value = 20000
$ws.Cells.Item(100, 16).Value = 20000
$ws.Cells.Item(100, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(100, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(100, 19).Value = 1000
$ws.Cells.Item(100, 20).Value = 20
